$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Copy the header formatting from D1 onto the new E1 header cell,
    # then set its text.
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial(-4122)
    $ws.Range("E1").Value = "%-age"

    # E2:E9 hold each category's share of the total points (C-column count
    # divided by the grand total in C10).
    $ws.Range("E2:E9").Formula = "=C2/C`$10"

    # E10 is the overall total, i.e. the sum of the column above it.
    $ws.Range("E10").Formula = "=SUM(E2:E9)"
}

$excel.CutCopyMode = 0
